$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Initial Data")
$ws.Range("M10").Value = "test"
